# Regenerate orders with updated distance/size codes.
# The underlying values follow a simple token substitution:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# These tokens appear (as substrings) inside the Condition, Filename_Left,
# Filename_Right, Distance and Size columns. Using Cells.Replace performs
# an exact, in-place substring substitution across every cell on the sheet,
# which reproduces the shared-string-table edit shown in the diff without
# disturbing any other content, formatting, or cell structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPart (match part of cell contents) = 2, xlByRows = 1
$xlPart = 2
$xlByRows = 1

$ws.Cells.Replace("D51", "D55", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("D64", "D69", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("D80", "D86", $xlPart, $xlByRows, $false, $false, $false, $false)
$ws.Cells.Replace("S30", "S31", $xlPart, $xlByRows, $false, $false, $false, $false)
